{"js": "// Apply the \"King of Africa\" copy refresh via Office.js (Word JS API).\n// Each entry is an exact-text replacement; matchCase keeps the search\n// precise (the body contains several similarly-worded list items).\nconst replacements = [\n  [\n    \"Play King of Africa Slot Game for Free - Review\",\n    \"Play King of Africa - Free Slot Game\",\n  ],\n  [\n    \"Exciting bonus mode triggered by finding three or more Scatters\",\n    \"Savannah theme with vibrant colors\",\n  ],\n  [\n    \"Engaging respin feature triggered by finding two Scatters\",\n    \"Exciting bonus features\",\n  ],\n  [\n    \"Savannah-themed design and vibrant color scheme add to the game's appeal\",\n    \"Wide range of symbols and combinations\",\n  ],\n  [\n    \"Variety of regular symbols keep gameplay interesting\",\n    \"Plenty of other Africa-themed slots to explore\",\n  ],\n  [\n    \"Relatively low RTP compared to the market's current standards\",\n    \"Low RTP compared to industry standards\",\n  ],\n  [\n    \"Lack of constant background music may make the game less immersive\",\n    \"Lack of constant background music\",\n  ],\n  [\n    \"Explore the pros and cons of King of Africa, a 5\\u00d73 online slot game with 20 paylines and exciting bonus features. Play for free and experience the savannah-themed design.\",\n    \"Read our review of King of Africa, an exciting slot game with vibrant colors and exciting bonus features. Play for free now!\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"Play King of Africa Slot Game for Free - Review\", \"Play King of Africa - Free Slot Game\"),\n  @(\"Exciting bonus mode triggered by finding three or more Scatters\", \"Savannah theme with vibrant colors\"),\n  @(\"Engaging respin feature triggered by finding two Scatters\", \"Exciting bonus features\"),\n  @(\"Savannah-themed design and vibrant color scheme add to the game's appeal\", \"Wide range of symbols and combinations\"),\n  @(\"Variety of regular symbols keep gameplay interesting\", \"Plenty of other Africa-themed slots to explore\"),\n  @(\"Relatively low RTP compared to the market's current standards\", \"Low RTP compared to industry standards\"),\n  @(\"Lack of constant background music may make the game less immersive\", \"Lack of constant background music\"),\n  @(\"Explore the pros and cons of King of Africa, a 5\u00d73 online slot game with 20 paylines and exciting bonus features. Play for free and experience the savannah-themed design.\", \"Read our review of King of Africa, an exciting slot game with vibrant colors and exciting bonus features. Play for free now!\")\n)\n\nforeach ($pair in $replacements) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $new\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.MatchCase = $true\n  $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
